{"js": "// The paragraph in question originally holds two runs split around a\n// \"_GoBack\" bookmark:\n//   run1: \"\u591a\u4e91\u8f6c\u5c0f\u96e8\"\n//   <bookmarkStart/><bookmarkEnd name=\"_GoBack\"/>\n//   run2: \"\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\"\n//\n// The edit appends \"\u54e6\" to run1's text (merging run2's content into it)\n// and removes run2 entirely, while leaving the bookmark untouched.\n\nconst body = context.document.body;\n\n// 1) Find the second run's text and delete that range (leaving the\n//    bookmark, which sits immediately before it, in place).\nconst secondRunText = \"\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\";\nconst secondHits = body.search(secondRunText, { matchCase: true });\nsecondHits.load(\"items\");\nawait context.sync();\n\nif (secondHits.items.length === 0) {\n  throw new Error(\"Could not find the second run's text to remove\");\n}\nsecondHits.items[0].delete();\nawait context.sync();\n\n// 2) Find the first run's text (\"\u591a\u4e91\u8f6c\u5c0f\u96e8\") and append the merged\n//    continuation (with the new trailing \"\u54e6\") right after it.\nconst firstRunText = \"\u591a\u4e91\u8f6c\u5c0f\u96e8\";\nconst firstHits = body.search(firstRunText, { matchCase: true });\nfirstHits.load(\"items\");\nawait context.sync();\n\nif (firstHits.items.length === 0) {\n  throw new Error(\"Could not find the first run's text\");\n}\nfirstHits.items[0].insertText(\n  \"\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\u54e6\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "# The target paragraph has two runs split around a \"_GoBack\" bookmark:\n#   run1: \"\u591a\u4e91\u8f6c\u5c0f\u96e8\"\n#   <bookmarkStart/><bookmarkEnd w:name=\"_GoBack\"/>\n#   run2: \"\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\"\n#\n# The edit merges run2's text into run1 (appending a trailing \"\u54e6\") and\n# removes run2 entirely, leaving the bookmark in place.\n\n$d = $word.ActiveDocument\n\n$secondRunText = \"\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\"\n$firstRunText  = \"\u591a\u4e91\u8f6c\u5c0f\u96e8\"\n$appendedText  = \"\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\u54e6\"\n\n# 1) Find the second run's text and delete it outright (the bookmark sits\n#    right before it and is left untouched).\n$delRange = $d.Content\n$found = $delRange.Find.Execute($secondRunText)\nif ($found) {\n    $delRange.Delete()\n}\n\n# 2) Find the first run's text and append the merged continuation\n#    (including the new trailing \"\u54e6\") right after it.\n$insRange = $d.Content\n$found2 = $insRange.Find.Execute($firstRunText)\nif ($found2) {\n    $insRange.Collapse(0)\n    $insRange.InsertAfter($appendedText)\n}\n"}
